# Update workbook version strings from
#   "mines - January 30 (built on February 02 2026 12.49.33 EST)"
# to
#   "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

# --- Sheet "About" ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: $newVersion"

$wsAbout.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Zhaozhuang Coal Mine, China, M0430, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- Sheet "Boundaries and methane sources" ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

$lastRow = $wsData.Cells.Item($wsData.Rows.Count, 19).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $wsData.Cells.Item($r, 19).Value = $newVersion
}
